# Progress on merging P&L programmes
# Consolidate the portfolio table: refresh prices, swap "USA" for "United
# States", add Country/Industry/Qty. shares columns for every holding, and
# drop the old ticker-only scratch rows that used to trail the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Apple) : refresh current price, fix country label ----------
$ws.Range("C2").Value = 146.21
$ws.Range("F2").Value = "United States"

# --- Row 3 (Amazon) -----------------------------------------------------
$ws.Range("C3").Value = 94.185000000000002
$ws.Range("F3").Value = "United States"
$ws.Range("G3").Value = "Consumer Cyclical"
$ws.Range("H3").Value = 8

# --- Row 4 (Alphabet) ----------------------------------------------------
$ws.Range("C4").Value = 100.18
$ws.Range("F4").Value = "United States"
$ws.Range("G4").Value = "Communication Services"
$ws.Range("H4").Value = 15

# --- Row 5 (NVIDIA) -------------------------------------------------------
$ws.Range("C5").Value = 165.63499999999999
$ws.Range("F5").Value = "United States"
$ws.Range("G5").Value = "Technology"
$ws.Range("H5").Value = 5

# --- Row 6 (Intel) --------------------------------------------------------
$ws.Range("C6").Value = 28.864999999999998
$ws.Range("F6").Value = "United States"
$ws.Range("G6").Value = "Technology"
$ws.Range("H6").Value = 2

# --- Row 7 (Pfizer) -------------------------------------------------------
$ws.Range("C7").Value = 51.02
$ws.Range("F7").Value = "United States"
$ws.Range("G7").Value = "Healthcare"
$ws.Range("H7").Value = 1

# --- Row 8 : was a bare ticker row (ENPH), now a full Enphase holding ----
$ws.Range("A8").Value = "Enphase Energy, Inc."
$ws.Range("B8").Value = "ENPH"
$ws.Range("C8").Value = 332.63
$ws.Range("F8").Value = "United States"
$ws.Range("G8").Value = "Technology"
$ws.Range("H8").Value = 12

# --- Row 9 : was a bare ticker row (CSCO), now a full Cisco holding ------
$ws.Range("A9").Value = "Cisco Systems, Inc."
$ws.Range("B9").Value = "CSCO"
$ws.Range("C9").Value = 49.365000000000002
$ws.Range("F9").Value = "United States"
$ws.Range("G9").Value = "Technology"
$ws.Range("H9").Value = 0

# --- Row 10 : was a bare ticker row (V), now a full Visa holding ---------
$ws.Range("A10").Value = "Visa Inc."
$ws.Range("B10").Value = "V"
$ws.Range("C10").Value = 216.01
$ws.Range("F10").Value = "United States"
$ws.Range("G10").Value = "Financial Services"
$ws.Range("H10").Value = 19

# --- Drop the old ticker-only scratch rows (11-22: QCOM, TSLA, LLOY.L, ---
# --- SHEL, AZN, ULVR, HSBA, BP, LSEG, JPM, UBSG, MSFT) --------------------
$ws.Range("A11:K22").ClearContents()

# --- Column widths for the new Country / Industry columns ----------------
$ws.Columns.Item(6).ColumnWidth = 17.9
$ws.Columns.Item(7).ColumnWidth = 21.8

# --- Selection / view housekeeping ---------------------------------------
$ws.Range("E1").Select()
